# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Column G held the "aragon" dimension/mapping. It is being replaced with the
# new curated "refArea" dimension, which uses a URI (not a skos:Concept) and
# has no separate mapping workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: iaest-dimension:aragon -> sdmx-dimension:refArea
$ws.Range("G2").Value = "sdmx-dimension:refArea"

# Row 4: skos:Concept -> URI-Comunidad (new value type for the refArea dimension)
$ws.Range("G4").Value = "URI-Comunidad"

# Row 5: drop the mapping workbook reference (mapping-aragon.xlsx) entirely,
# leaving the cell empty like the other non-mapped columns (e.g. D5).
$ws.Range("G5").Clear()

$wb.Save()
